# Apply updated crypto price/volume figures (Coinranking snapshot refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.030.37'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '2.300.94'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'300.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = "'97.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = "'0.507"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = "'33.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.18%  '
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").Value = "'49.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.11%  '
$ws.Range("E13").Value = '  +3.12%  '
$ws.Range("D14").Value = "'17.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +13.02%  '
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '2.654.55'
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("D17").Value = '2.301.79'
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").Value = "'0.808"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.02%  '
$ws.Range("D19").Value = '42.983.92'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = "'11.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").Value = "'6.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("D23").Value = "'67.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("D24").Value = "'236.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +6.03%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").Value = "'24.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D30").Value = "'166.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").Value = "'33.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("D34").Value = "'4.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("D35").Value = "'4.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.96%  '
$ws.Range("D36").Value = "'2.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = "'16.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.25%  '
$ws.Range("D38").Value = "'0.0697"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("E39").Value = '  +0.63%  '
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").Value = "'2.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").Value = '1.982.10'
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").Value = "'0.0283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("D46").Value = "'9.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").Value = "'17.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("D49").Value = '2.534.34'
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").Value = "'53.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("D51").Value = "'4.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.11%  '
